# Enhance config file handling
# Append a new daily log entry (row 61) to each of the 4 worksheets,
# duplicating the last existing row (60) but advancing the timestamp.

$wb = $excel.ActiveWorkbook

$sheetNames = @("FE_LFT_#1", "FE_LFT_#2", "FE_PLT_#1", "FE_PLT_#2")

$newDate = 45847.495

$rowData = @{
    "FE_LFT_#1" = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x50"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 336
        I = 15
    }
    "FE_LFT_#2" = @{
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x60"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 352
        I = 14
    }
    "FE_PLT_#1" = @{
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x67"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 103
        I = 3
    }
    "FE_PLT_#2" = @{
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x67"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 103
        I = 3
    }
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $data = $rowData[$name]

    $ws.Cells.Item(61, 1).Value = $newDate
    $ws.Cells.Item(61, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item(61, 2).Value = $data.B
    $ws.Cells.Item(61, 3).Value = $data.C
    $ws.Cells.Item(61, 4).Value = $data.D
    $ws.Cells.Item(61, 5).Value = $data.E
    $ws.Cells.Item(61, 6).Value = $data.F
    $ws.Cells.Item(61, 7).Value = $data.G
    $ws.Cells.Item(61, 8).Value = $data.H
    $ws.Cells.Item(61, 9).Value = $data.I
}
